# Atualização de bases das ligas, do dia: 22-05-2024 às 20:16
# Swap the data of three pairs of rows (everything except column A, which
# just holds the running row index) back to their correct order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData($ws, $row1, $row2, $firstCol, $lastCol) {
    # Read all values from both rows first (columns firstCol..lastCol)
    $values1 = @()
    $values2 = @()
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $values1 += , $ws.Cells.Item($row1, $c).Value2
        $values2 += , $ws.Cells.Item($row2, $c).Value2
    }
    # Write them back swapped
    $i = 0
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws.Cells.Item($row1, $c).Value = $values2[$i]
        $ws.Cells.Item($row2, $c).Value = $values1[$i]
        $i++
    }
}

# Columns B (2) through AB (28) are swapped; column A (1) keeps the
# original sequential row number for each physical row.
Swap-RowData $ws 74 75 2 28
Swap-RowData $ws 140 141 2 28
Swap-RowData $ws 142 143 2 28
